# ============================================================================
# Add 2022-Q3 data:
#  1. Duplicate the "2022-Q2" worksheet, place the copy right before it, and
#     rename the copy to "2022-Q3" (all later quarter sheets keep their names,
#     they simply shift right in the tab order).
#  2. Overwrite the new "2022-Q3" sheet's holdings table (rows 2-22) with the
#     Q3 fund figures.
#  3. Prepend a "2022-Q3" row to the "总计" summary sheet and rewrite the
#     index column + the rows below so everything keeps lining up.
# ============================================================================

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "2022-Q2" -> "2022-Q3" --------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# --- 2. Rewrite the "2022-Q3" holdings table -----------------------------------
# Columns B,C,D,E,F,G are stored as text in this workbook (fund code/name/size/
# position/weight/value); a leading apostrophe keeps numeric-looking text (fund
# codes, "44.00", ...) from being auto-converted to a number and losing leading
# zeros / trailing zeros. Column H (rank) is a genuine number.
$q3.Range("B2").Value = "'513060"
$q3.Range("C2").Value = "博时恒生医疗保健ETF（QDII）"
$q3.Range("D2").Value = "'44.00"
$q3.Range("E2").Value = "'99.48"
$q3.Range("F2").Value = "'4.60"
$q3.Range("G2").Value = "'2.0240"
$q3.Range("H2").Value = 5
$q3.Range("B3").Value = "'002121"
$q3.Range("C3").Value = "广发沪港深新起点股票A"
$q3.Range("D3").Value = "'25.97"
$q3.Range("E3").Value = "'85.60"
$q3.Range("F3").Value = "'4.56"
$q3.Range("G3").Value = "'1.1842"
$q3.Range("H3").Value = 7
$q3.Range("B4").Value = "'011338"
$q3.Range("C4").Value = "兴全合远两年持有期混合A"
$q3.Range("D4").Value = "'31.13"
$q3.Range("E4").Value = "'84.61"
$q3.Range("F4").Value = "'2.97"
$q3.Range("G4").Value = "'0.9246"
$q3.Range("H4").Value = 9
$q3.Range("B5").Value = "'009993"
$q3.Range("C5").Value = "嘉实前沿创新混合"
$q3.Range("D5").Value = "'14.33"
$q3.Range("E5").Value = "'89.66"
$q3.Range("F5").Value = "'4.02"
$q3.Range("G5").Value = "'0.5761"
$q3.Range("H5").Value = 8
$q3.Range("B6").Value = "'010387"
$q3.Range("C6").Value = "易方达医药生物股票A"
$q3.Range("D6").Value = "'16.27"
$q3.Range("E6").Value = "'84.59"
$q3.Range("F6").Value = "'3.48"
$q3.Range("G6").Value = "'0.5662"
$q3.Range("H6").Value = 8
$q3.Range("B7").Value = "'010388"
$q3.Range("C7").Value = "易方达医药生物股票C"
$q3.Range("D7").Value = "'4.82"
$q3.Range("E7").Value = "'84.59"
$q3.Range("F7").Value = "'3.48"
$q3.Range("G7").Value = "'0.1677"
$q3.Range("H7").Value = 8
$q3.Range("B8").Value = "'012826"
$q3.Range("C8").Value = "工银聚宁9个月持有期混合A"
$q3.Range("D8").Value = "'8.98"
$q3.Range("E8").Value = "'26.05"
$q3.Range("F8").Value = "'1.34"
$q3.Range("G8").Value = "'0.1203"
$q3.Range("H8").Value = 7
$q3.Range("B9").Value = "'513700"
$q3.Range("C9").Value = "鹏华中证港股通医药卫生综合ETF"
$q3.Range("D9").Value = "'2.96"
$q3.Range("E9").Value = "'94.14"
$q3.Range("F9").Value = "'3.79"
$q3.Range("G9").Value = "'0.1122"
$q3.Range("H9").Value = 6
$q3.Range("B10").Value = "'007718"
$q3.Range("C10").Value = "中银创新医疗混合A"
$q3.Range("D10").Value = "'2.82"
$q3.Range("E10").Value = "'92.35"
$q3.Range("F10").Value = "'3.74"
$q3.Range("G10").Value = "'0.1055"
$q3.Range("H10").Value = 8
$q3.Range("B11").Value = "'159892"
$q3.Range("C11").Value = "华夏恒生香港上市生物科技ETF（QDII）"
$q3.Range("D11").Value = "'1.58"
$q3.Range("E11").Value = "'99.13"
$q3.Range("F11").Value = "'5.05"
$q3.Range("G11").Value = "'0.0798"
$q3.Range("H11").Value = 5
$q3.Range("B12").Value = "'513120"
$q3.Range("C12").Value = "广发中证香港创新药（QDII-ETF）"
$q3.Range("D12").Value = "'1.09"
$q3.Range("E12").Value = "'98.58"
$q3.Range("F12").Value = "'7.13"
$q3.Range("G12").Value = "'0.0777"
$q3.Range("H12").Value = 5
$q3.Range("B13").Value = "'513280"
$q3.Range("C13").Value = "汇添富恒生香港上市生物科技ETF（QDII）"
$q3.Range("D13").Value = "'1.51"
$q3.Range("E13").Value = "'100.14"
$q3.Range("F13").Value = "'5.11"
$q3.Range("G13").Value = "'0.0772"
$q3.Range("H13").Value = 5
$q3.Range("B14").Value = "'470888"
$q3.Range("C14").Value = "汇添富香港优势精选混合（QDII）"
$q3.Range("D14").Value = "'1.63"
$q3.Range("E14").Value = "'78.50"
$q3.Range("F14").Value = "'3.44"
$q3.Range("G14").Value = "'0.0561"
$q3.Range("H14").Value = 10
$q3.Range("B15").Value = "'011339"
$q3.Range("C15").Value = "兴全合远两年持有期混合C"
$q3.Range("D15").Value = "'1.46"
$q3.Range("E15").Value = "'84.61"
$q3.Range("F15").Value = "'2.97"
$q3.Range("G15").Value = "'0.0434"
$q3.Range("H15").Value = 9
$q3.Range("B16").Value = "'513200"
$q3.Range("C16").Value = "易方达中证港股通医药卫生综合ETF"
$q3.Range("D16").Value = "'0.77"
$q3.Range("E16").Value = "'95.67"
$q3.Range("F16").Value = "'3.98"
$q3.Range("G16").Value = "'0.0306"
$q3.Range("H16").Value = 6
$q3.Range("B17").Value = "'010024"
$q3.Range("C17").Value = "广发沪港深新起点股票C"
$q3.Range("D17").Value = "'0.49"
$q3.Range("E17").Value = "'85.60"
$q3.Range("F17").Value = "'4.56"
$q3.Range("G17").Value = "'0.0223"
$q3.Range("H17").Value = 7
$q3.Range("B18").Value = "'159776"
$q3.Range("C18").Value = "银华中证港股通医药卫生综合ETF"
$q3.Range("D18").Value = "'0.52"
$q3.Range("E18").Value = "'92.74"
$q3.Range("F18").Value = "'3.74"
$q3.Range("G18").Value = "'0.0194"
$q3.Range("H18").Value = 6
$q3.Range("B19").Value = "'159718"
$q3.Range("C19").Value = "平安中证港股通医药卫生综合ETF"
$q3.Range("D19").Value = "'0.53"
$q3.Range("E19").Value = "'90.14"
$q3.Range("F19").Value = "'3.63"
$q3.Range("G19").Value = "'0.0192"
$q3.Range("H19").Value = 6
$q3.Range("B20").Value = "'012827"
$q3.Range("C20").Value = "工银聚宁9个月持有期混合C"
$q3.Range("D20").Value = "'0.64"
$q3.Range("E20").Value = "'26.05"
$q3.Range("F20").Value = "'1.34"
$q3.Range("G20").Value = "'0.0086"
$q3.Range("H20").Value = 7
$q3.Range("B21").Value = "'010500"
$q3.Range("C21").Value = "中银创新医疗混合C"
$q3.Range("D21").Value = "'0.18"
$q3.Range("E21").Value = "'92.35"
$q3.Range("F21").Value = "'3.74"
$q3.Range("G21").Value = "'0.0067"
$q3.Range("H21").Value = 8
$q3.Range("B22").Value = "'006603"
$q3.Range("C22").Value = "嘉实互融精选股票"
$q3.Range("D22").Value = "'0.12"
$q3.Range("E22").Value = "'82.85"
$q3.Range("F22").Value = "'5.21"
$q3.Range("G22").Value = "'0.0063"
$q3.Range("H22").Value = 2

# --- 3. Prepend the 2022-Q3 total row on "总计" --------------------------------
$sum = $wb.Worksheets.Item("总计")

# Row 9 is brand new - copy the bordered/centered style used by the other index
# cells in column A onto it before the values are written.
$sum.Range("A2").Copy()
$sum.Range("A9").PasteSpecial(-4122)

$sum.Range("A2").Value = 0
$sum.Range("B2").Value = "2022-Q3"
$sum.Range("C2").Value = 21
$sum.Range("D2").Value = 6.23
$sum.Range("A3").Value = 1
$sum.Range("B3").Value = "2022-Q2"
$sum.Range("C3").Value = 21
$sum.Range("D3").Value = 13.54
$sum.Range("A4").Value = 2
$sum.Range("B4").Value = "2022-Q1"
$sum.Range("C4").Value = 16
$sum.Range("D4").Value = 4.94
$sum.Range("A5").Value = 3
$sum.Range("B5").Value = "2021-Q4"
$sum.Range("C5").Value = 12
$sum.Range("D5").Value = 11.52
$sum.Range("A6").Value = 4
$sum.Range("B6").Value = "2021-Q3"
$sum.Range("C6").Value = 38
$sum.Range("D6").Value = 32.96
$sum.Range("A7").Value = 5
$sum.Range("B7").Value = "2021-Q2"
$sum.Range("C7").Value = 51
$sum.Range("D7").Value = 33.61
$sum.Range("A8").Value = 6
$sum.Range("B8").Value = "2021-Q1"
$sum.Range("C8").Value = 49
$sum.Range("D8").Value = 34.59
$sum.Range("A9").Value = 7
$sum.Range("B9").Value = "2020-Q4"
$sum.Range("C9").Value = 30
$sum.Range("D9").Value = 12.77

# Restore "总计" as the active sheet (it was active before the edit).
$sum.Select()
